# Commit: "All results categorized and cleaned up"
#
# Sheet1 had a trailing "averages" block (row 117: per-column AVERAGE()s,
# row 119: a SUM-of-averages "total average time:" cell). That block is
# pulled out into its own "averages" worksheet, laid out horizontally as
# a one-row table (with the old "total average time:" label renamed to
# "total_time" and reused as that column's header), and the now-empty
# rows are removed from Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the shared string in place (D119, "total average time:") to
# "total_time" *before* copying it into the new sheet's header and
# *before* deleting the row it currently lives in, so the same shared
# string slot is reused rather than dropped/recreated.
$ws1.Range("D119").Value = "total_time"

# New worksheet, positioned right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "averages"

# Header row: reuse Sheet1's existing column headers (E1:H1, J1:K1) plus
# the renamed "total_time" label, in the new column order
# A..D = E..H (the four AVERAGE() columns), E = total_time, F..G = J..K.
$ws2.Range("A1").Value2 = $ws1.Range("E1").Value2
$ws2.Range("B1").Value2 = $ws1.Range("F1").Value2
$ws2.Range("C1").Value2 = $ws1.Range("G1").Value2
$ws2.Range("D1").Value2 = $ws1.Range("H1").Value2
$ws2.Range("E1").Value2 = $ws1.Range("D119").Value2
$ws2.Range("F1").Value2 = $ws1.Range("J1").Value2
$ws2.Range("G1").Value2 = $ws1.Range("K1").Value2

# Data row: same AVERAGE()/SUM() formulas that used to live in Sheet1
# rows 117/119, now pointing at Sheet1 explicitly and laid out across a
# single row.
$ws2.Range("A2").Formula = "=AVERAGE(Sheet1!E2:E116)"
$ws2.Range("B2").Formula = "=AVERAGE(Sheet1!F2:F116)"
$ws2.Range("C2").Formula = "=AVERAGE(Sheet1!G2:G116)"
$ws2.Range("D2").Formula = "=AVERAGE(Sheet1!H2:H116)"
$ws2.Range("E2").Formula = "=SUM(A2:D2)"
$ws2.Range("F2").Formula = "=AVERAGE(Sheet1!J2:J116)"
$ws2.Range("G2").Formula = "=AVERAGE(Sheet1!K2:K116)"

# Cosmetic column widths to roughly match the authored sheet (best effort).
$ws2.Columns.Item(1).ColumnWidth = 22.59
$ws2.Columns.Item(2).ColumnWidth = 24.59
$ws2.Columns.Item(3).ColumnWidth = 23.75
$ws2.Columns.Item(4).ColumnWidth = 21.75
$ws2.Columns.Item(5).ColumnWidth = 16.59

# Now that the figures live on "averages", drop the old summary rows
# from Sheet1 (bottom-up, so the row numbers don't shift underneath us).
$ws1.Range("A119").EntireRow.Delete()
$ws1.Range("A117").EntireRow.Delete()

# Selection / active-sheet bookkeeping to match the saved file: land on
# "averages" (now the active tab) at E3, leaving Sheet1's selection at E2.
[void]$ws1.Range("E2").Select()
[void]$ws2.Range("E3").Select()
